# Junction_Flooding_144.xlsx edit
# - Round row 5 flow values to 2 decimal places (custom accuracy)
# - Remove row 6 (data trimmed / re-sampled -> fewer rows)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 5 values to their 2-decimal-place rounded equivalents
$ws.Range("B5").Value = 3.72
$ws.Range("C5").Value = 2.79
$ws.Range("D5").Value = 0.31
$ws.Range("E5").Value = 8.12
$ws.Range("F5").Value = 6.81
$ws.Range("G5").Value = 3.09
$ws.Range("H5").Value = 13.87
$ws.Range("I5").Value = 4.48
$ws.Range("J5").Value = 2.26
$ws.Range("K5").Value = 2.93
$ws.Range("L5").Value = 3.19
$ws.Range("M5").Value = 3.71
$ws.Range("N5").Value = 1.35
$ws.Range("O5").Value = 2.89
$ws.Range("P5").Value = 4.31
$ws.Range("Q5").Value = 2.23
$ws.Range("R5").Value = 0.23
$ws.Range("S5").Value = 0.25
$ws.Range("T5").Value = 38.29
$ws.Range("U5").Value = 8.32
$ws.Range("V5").Value = 2.66
$ws.Range("W5").Value = 6
$ws.Range("X5").Value = 2.87
$ws.Range("Y5").Value = 0.41
$ws.Range("Z5").Value = 6.73
$ws.Range("AA5").Value = 2.43
$ws.Range("AB5").Value = 2.28
$ws.Range("AC5").Value = 2.61
$ws.Range("AD5").Value = 3.92
$ws.Range("AE5").Value = 0.52
$ws.Range("AF5").Value = 13
$ws.Range("AG5").Value = 1.58
$ws.Range("AH5").Value = 3.47

# Delete the now-obsolete last row of data (row 6), which also
# shrinks the sheet dimension from A1:AH6 to A1:AH5
$ws.Rows(6).Delete()
